$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026327809050455
$ws.Range("D2").Value = 1.033934170463273
$ws.Range("E2").Value = 1.026562891491226
$ws.Range("F2").Value = 1.041029263270843
$ws.Range("I2").Value = 1.030972423753723
$ws.Range("J2").Value = 1.031491865453045
$ws.Range("K2").Value = 1.036735092774435
$ws.Range("L2").Value = 1.02938517384718
$ws.Range("M2").Value = 1.043809936069994
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027284868655318
$ws.Range("D3").Value = 1.034652058386466
$ws.Range("E3").Value = 1.027375804996841
$ws.Range("F3").Value = 1.041943570300468
$ws.Range("I3").Value = 1.031104258443841
$ws.Range("J3").Value = 1.032088599008431
$ws.Range("K3").Value = 1.037262307209987
$ws.Range("L3").Value = 1.030005613419861
$ws.Range("M3").Value = 1.044534514629731
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027904694304471
$ws.Range("D4").Value = 1.035116844112886
$ws.Range("E4").Value = 1.027902660592929
$ws.Range("F4").Value = 1.042535956090852
$ws.Range("I4").Value = 1.031188277057087
$ws.Range("J4").Value = 1.032474676121379
$ws.Range("K4").Value = 1.037603044275085
$ws.Range("L4").Value = 1.030407276591723
$ws.Range("M4").Value = 1.045003503623353
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028165398288454
$ws.Range("D5").Value = 1.035312301825979
$ws.Range("E5").Value = 1.028124351883175
$ws.Range("F5").Value = 1.042785177462964
$ws.Range("I5").Value = 1.031223289755672
$ws.Range("J5").Value = 1.032636970101201
$ws.Range("K5").Value = 1.037746191859673
$ws.Range("L5").Value = 1.030576181941058
$ws.Range("M5").Value = 1.045200698402347
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028209179187396
$ws.Range("D6").Value = 1.035345123617063
$ws.Range("E6").Value = 1.028161586584383
$ws.Range("F6").Value = 1.04282703348028
$ws.Range("I6").Value = 1.031229150421788
$ws.Range("J6").Value = 1.032664219192287
$ws.Range("K6").Value = 1.03777022117406
$ws.Range("L6").Value = 1.030604544567486
$ws.Range("M6").Value = 1.045233810101291
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027908177337089
$ws.Range("D7").Value = 1.035119455587451
$ws.Range("E7").Value = 1.027905622053969
$ws.Range("F7").Value = 1.042539285484224
$ws.Range("I7").Value = 1.031188746112258
$ws.Range("J7").Value = 1.032476844753196
$ws.Range("K7").Value = 1.037604957407659
$ws.Range("L7").Value = 1.030409533333274
$ws.Range("M7").Value = 1.045006138426275
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026651138580372
$ws.Range("D8").Value = 1.034176728455598
$ws.Range("E8").Value = 1.026837443372791
$ws.Range("F8").Value = 1.041338097704045
$ws.Range("I8").Value = 1.031017244102342
$ws.Range("J8").Value = 1.031693543971737
$ws.Range("K8").Value = 1.03691335075156
$ws.Range("L8").Value = 1.029594812809339
$ws.Range("M8").Value = 1.044054781396931
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02444028580838
$ws.Range("D9").Value = 1.032517611180028
$ws.Range("E9").Value = 1.024961718060623
$ws.Range("F9").Value = 1.039227398532732
$ws.Range("I9").Value = 1.030705201477232
$ws.Range("J9").Value = 1.03031293423995
$ws.Range("K9").Value = 1.035691590962966
$ws.Range("L9").Value = 1.028160729648755
$ws.Range("M9").Value = 1.042379483940625
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.022969271401685
$ws.Range("D10").Value = 1.031413023755124
$ws.Range("E10").Value = 1.023715714828008
$ws.Range("F10").Value = 1.037824351427937
$ws.Range("I10").Value = 1.030490594060963
$ws.Range("J10").Value = 1.029392362998891
$ws.Range("K10").Value = 1.034875091508913
$ws.Range("L10").Value = 1.027205785615195
$ws.Range("M10").Value = 1.041263444447965
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022333002344295
$ws.Range("D11").Value = 1.030935096906643
$ws.Range("E11").Value = 1.023177261423274
$ws.Range("F11").Value = 1.037217803750861
$ws.Range("I11").Value = 1.030396113085392
$ws.Range("J11").Value = 1.028993718768648
$ws.Range("K11").Value = 1.034521079020396
$ws.Range("L11").Value = 1.026792561138262
$ws.Range("M11").Value = 1.040780398568353
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022096768213192
$ws.Range("D12").Value = 1.0307576301433
$ws.Range("E12").Value = 1.022977418603912
$ws.Range("F12").Value = 1.036992653754951
$ws.Range("I12").Value = 1.030360785616719
$ws.Range("J12").Value = 1.028845640889971
$ws.Range("K12").Value = 1.034389514542341
$ws.Range("L12").Value = 1.026639113164868
$ws.Range("M12").Value = 1.040601006099229
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022147436492011
$ws.Range("D13").Value = 1.030795694784428
$ws.Range("E13").Value = 1.023020278180429
$ws.Range("F13").Value = 1.037040942435956
$ws.Range("I13").Value = 1.03036837401521
$ws.Range("J13").Value = 1.028877404247788
$ws.Range("K13").Value = 1.034417738651237
$ws.Range("L13").Value = 1.026672026358896
$ws.Range("M13").Value = 1.040639484912228
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022313473019908
$ws.Range("D14").Value = 1.030920426286735
$ws.Range("E14").Value = 1.023160739025567
$ws.Range("F14").Value = 1.037199189737808
$ws.Range("I14").Value = 1.030393197656459
$ws.Range("J14").Value = 1.028981478674111
$ws.Range("K14").Value = 1.034510205250393
$ws.Range("L14").Value = 1.026779876218173
$ws.Range("M14").Value = 1.040765569265248
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02241578748183
$ws.Range("D15").Value = 1.03099728501819
$ws.Range("E15").Value = 1.023247303194682
$ws.Range("F15").Value = 1.037296710886119
$ws.Range("I15").Value = 1.030408461458444
$ws.Range("J15").Value = 1.029045601906752
$ws.Range("K15").Value = 1.034567167931664
$ws.Range("L15").Value = 1.026846331672059
$ws.Range("M15").Value = 1.040843258310764
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023011513093559
$ws.Range("D16").Value = 1.031444750037035
$ws.Range("E16").Value = 1.023751472955808
$ws.Range("F16").Value = 1.037864626790365
$ws.Range("I16").Value = 1.030496831744453
$ws.Range("J16").Value = 1.029418819136944
$ws.Range("K16").Value = 1.034898576490777
$ws.Range("L16").Value = 1.027233215794957
$ws.Range("M16").Value = 1.041295507070553
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.023385381352319
$ws.Range("D17").Value = 1.031725532240765
$ws.Range("E17").Value = 1.024068013886138
$ws.Range("F17").Value = 1.038221129082795
$ws.Range("I17").Value = 1.030551848264283
$ws.Range("J17").Value = 1.029652920879071
$ws.Range("K17").Value = 1.035106337156902
$ws.Range("L17").Value = 1.027475971700072
$ws.Range("M17").Value = 1.041579246953994
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023603518887261
$ws.Range("D18").Value = 1.031889343043635
$ws.Range("E18").Value = 1.024252750410497
$ws.Range("F18").Value = 1.038429165392148
$ws.Range("I18").Value = 1.03058378844283
$ws.Range("J18").Value = 1.029789465427152
$ws.Range("K18").Value = 1.035227475710364
$ws.Range("L18").Value = 1.02761759334359
$ws.Range("M18").Value = 1.041744767476727
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.023677909403234
$ws.Range("D19").Value = 1.031945204236763
$ws.Range("E19").Value = 1.024315758314633
$ws.Range("F19").Value = 1.038500116391747
$ws.Range("I19").Value = 1.030594653755203
$ws.Range("J19").Value = 1.029836023044832
$ws.Range("K19").Value = 1.03526877320774
$ws.Range("L19").Value = 1.027665887081755
$ws.Range("M19").Value = 1.041801209026575
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023345261941102
$ws.Range("D20").Value = 1.031695403306565
$ws.Range("E20").Value = 1.024034041318989
$ws.Range("F20").Value = 1.038182869970185
$ws.Range("I20").Value = 1.030545961026537
$ws.Range("J20").Value = 1.02962780428483
$ws.Range("K20").Value = 1.035084051022423
$ws.Range("L20").Value = 1.027449923568113
$ws.Range("M20").Value = 1.041548802287706
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022264576508262
$ws.Range("D21").Value = 1.030883694395593
$ws.Range("E21").Value = 1.02311937230395
$ws.Range("F21").Value = 1.037152585731975
$ws.Range("I21").Value = 1.030385894141299
$ws.Range("J21").Value = 1.028950831442865
$ws.Range("K21").Value = 1.03448297802342
$ws.Range("L21").Value = 1.02674811595503
$ws.Range("M21").Value = 1.040728439664847
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021585711882599
$ws.Range("D22").Value = 1.030373668667638
$ws.Range("E22").Value = 1.022545225884821
$ws.Range("F22").Value = 1.036505666866172
$ws.Range("I22").Value = 1.030283905772786
$ws.Range("J22").Value = 1.028525171040694
$ws.Range("K22").Value = 1.034104663907472
$ws.Range("L22").Value = 1.026307105147277
$ws.Range("M22").Value = 1.040212832490367
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.021945533204537
$ws.Range("D23").Value = 1.030644011395868
$ws.Range("E23").Value = 1.022849501991879
$ws.Range("F23").Value = 1.036848528635186
$ws.Range("I23").Value = 1.030338099332293
$ws.Range("J23").Value = 1.028750823313415
$ws.Range("K23").Value = 1.0343052525292
$ws.Range("L23").Value = 1.026540869892108
$ws.Range("M23").Value = 1.040486147434354
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023363389970664
$ws.Range("D24").Value = 1.031709017166144
$ws.Range("E24").Value = 1.024049391739925
$ws.Range("F24").Value = 1.038200157324313
$ws.Range("I24").Value = 1.030548621679468
$ws.Range("J24").Value = 1.029639153402105
$ws.Range("K24").Value = 1.035094121305599
$ws.Range("L24").Value = 1.027461693516766
$ws.Range("M24").Value = 1.041562558860712
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025011339488398
$ws.Range("D25").Value = 1.032946276347522
$ws.Range("E25").Value = 1.025445854320966
$ws.Range("F25").Value = 1.039772351116416
$ws.Range("I25").Value = 1.030787033984336
$ws.Range("J25").Value = 1.030669888353762
$ws.Range("K25").Value = 1.036007801163983
$ws.Range("L25").Value = 1.028531283282607
$ws.Range("M25").Value = 1.042812448537926